$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.085.20'
$ws.Range('E2').Value = '  -2.36%  '
$ws.Range('D3').Value = '1.822.03'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -1.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4223'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3677'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07213'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8402'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.78'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.03%  '
$ws.Range('D12').Value = '1.814.73'
$ws.Range('E12').Value = '  -1.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.650'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07067'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.282'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '89.47'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('E17').Value = '  -1.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008788'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.47%  '
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.93'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.77%  '
$ws.Range('D21').Value = '27.059.63'
$ws.Range('E21').Value = '  -2.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.120'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.84'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.54%  '
$ws.Range('D24').Value = '2.038.73'
$ws.Range('E24').Value = '  -1.94%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.976'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.97%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.63'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.226'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.25'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.232'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08748'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.25%  '
$ws.Range('E32').Value = '  -4.71%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7389'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.15%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.945'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.417'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9998'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.090'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.96%  '
$ws.Range('E38').Value = '  -1.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05242'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.319'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.876'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1686'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5027'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.590'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.55'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '106.25'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4714'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9996'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.29%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06350'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.889'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.646'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.00%  '
